$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.628.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.595.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.611.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.588.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.281.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +10.72%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.599"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.730.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
